$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Juan"
$ws.Range("B4").Value = "Emp"
$ws.Range("C4").Value = "Empleado"

$ws.Range("A5").Value = "Juan"
$ws.Range("B5").Value = "Cli"
$ws.Range("C5").Value = "Cliente"

$ws.Range("E7").Select()
